# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Change cell B11 on the "Rules" sheet from the text "R40" to the text "1".
#
# A plain `Range.Value = "1"` assignment gets auto-coerced to the number 1
# (same as typing "1" straight into a General-formatted cell in Excel), but
# the target cell must stay a *text* string (t="s" in the OOXML) with its
# existing cell style untouched. Writing a text-producing formula and then
# collapsing it down to a static value via Copy/PasteSpecial(values) gets
# us a genuine text cell without disturbing the cell's style index.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
$cell.Formula = "=""1"""
$cell.Copy()
$cell.PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false
